$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

$ws.Range("F8").Value = 'Both peers send strings (previously "boxes") - every other frame supports this basic function. `string` is restricted to the base "restricted string" codepoints.'
$ws.Range("F9").Value = 'Either peer can reset if they''ve given up on this stream. `reasonString` is restricted to the base "restricted string" codepoints.'

$ws.Range("F9").Select() | Out-Null
